$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.954.31"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "1.639.29"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -1.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.49"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("E6").Value = "  +1.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.66%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2561"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06364"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.51"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07758"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.280"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("D13").Value = "1.642.14"
$ws.Range("E13").Value = "  -0.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5436"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("D15").Value = "0.0₅7762"
$ws.Range("E15").Value = "  -1.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.27"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.26%  "
$ws.Range("D17").Value = "25.956.53"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "196.35"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.426"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.925"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("E22").Value = "  +1.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.879"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "140.97"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1193"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +5.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.845"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.65%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.234"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.04935"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.246"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.530"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.363"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.95%  "
$ws.Range("E35").Value = "  +0.69%  "
$ws.Range("D36").Value = "1.150.64"
$ws.Range("E36").Value = "  +2.03%  "
$ws.Range("E37").Value = "  -1.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5414"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01553"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.000"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.521"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.60%  "
$ws.Range("B42").Value = "BabyDogeCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D42").Value = "0.0₈126"
$ws.Range("E42").Value = "  +7.89%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8094"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.450"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "98.85"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.52%  "
$ws.Range("D46").Value = "1.776.26"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4521"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9992"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.83"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05050"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.000"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.00%  "
